# Auto-generated Excel COM-interop edit script
# Updates cached market-price figures on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled price-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 714728.4399999999
$ws.Range("I6").Value = 833350
$ws.Range("J6").Value = 2999
$ws.Range("K6").Value = 2500050
$ws.Range("L6").Value = 8997
$ws.Range("M6").Value = -2499938
$ws.Range("N6").Value = -9221

$ws.Range("H32").Value = 1500
$ws.Range("I32").Value = 1500
$ws.Range("K32").Value = 1500
$ws.Range("M32").Value = -1174

$ws.Range("H50").Value = 997
$ws.Range("J50").Value = 997
$ws.Range("L50").Value = 2991
$ws.Range("N50").Value = -3941

$ws.Range("H51").Value = 7819.5
$ws.Range("I51").Value = 7389
$ws.Range("J51").Value = 8250
$ws.Range("K51").Value = 7389
$ws.Range("L51").Value = 8250
$ws.Range("M51").Value = -6905
$ws.Range("N51").Value = -9218

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6818.6665
$ws.Range("I61").Value = 4313.778
$ws.Range("K61").Value = 4313.778
$ws.Range("M61").Value = -4101.778

$ws.Range("H92").Value = 86000
$ws.Range("I92").Value = 63000
$ws.Range("J92").Value = 97500
$ws.Range("K92").Value = 63000
$ws.Range("L92").Value = 97500
$ws.Range("M92").Value = -60504
$ws.Range("N92").Value = -102492

$ws.Range("H97").Value = 2007
$ws.Range("I97").Value = 2005
$ws.Range("K97").Value = 2005
$ws.Range("M97").Value = -1509

$ws.Range("H98").Value = 58138.75
$ws.Range("J98").Value = 58138.75
$ws.Range("L98").Value = 58138.75
$ws.Range("N98").Value = -64128.75

$ws.Range("H122").Value = 2552.2727
$ws.Range("I122").Value = 2258.5
$ws.Range("K122").Value = 6775.5
$ws.Range("M122").Value = -4325.5

$ws.Range("H136").Value = 6818.6665
$ws.Range("I136").Value = 4313.778
$ws.Range("K136").Value = 12941.334
$ws.Range("M136").Value = -10391.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 6424.75
$ws.Range("J100").Value = 6424.75
$ws.Range("L100").Value = 6424.75
$ws.Range("N100").Value = -8588.75

$ws.Range("H105").Value = 1900
$ws.Range("I105").Value = 1833.3334
$ws.Range("J105").Value = 2100
$ws.Range("K105").Value = 1833.3334
$ws.Range("L105").Value = 2100
$ws.Range("M105").Value = -86.33339999999998
$ws.Range("N105").Value = -5594

$ws.Range("H134").Value = 7687
$ws.Range("I134").Value = 1410.1428
$ws.Range("K134").Value = 4230.428400000001
$ws.Range("M134").Value = -1695.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1158.5
$ws.Range("I16").Value = 1179.2858
$ws.Range("K16").Value = 1179.2858
$ws.Range("M16").Value = -892.2858000000001

$ws.Range("H22").Value = 732
$ws.Range("I22").Value = 732
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 732
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -382
$ws.Range("N22").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H56").Value = 4296.5
$ws.Range("I56").Value = 4296.5
$ws.Range("K56").Value = 4296.5
$ws.Range("M56").Value = -3451.5

$ws.Range("H58").Value = 12142
$ws.Range("I58").Value = 3999
$ws.Range("K58").Value = 3999
$ws.Range("M58").Value = -3796

$ws.Range("H62").Value = 4997.5
$ws.Range("J62").Value = 5006
$ws.Range("L62").Value = 5006
$ws.Range("N62").Value = -6254

$ws.Range("H65").Value = 4997.5
$ws.Range("J65").Value = 5006
$ws.Range("L65").Value = 25030
$ws.Range("N65").Value = -31270

$ws.Range("H68").Value = 49799.25
$ws.Range("J68").Value = 49799.25
$ws.Range("L68").Value = 49799.25
$ws.Range("N68").Value = -51297.25

$ws.Range("H71").Value = 49799.25
$ws.Range("J71").Value = 49799.25
$ws.Range("L71").Value = 149397.75
$ws.Range("N71").Value = -156885.75

$ws.Range("H107").Value = 595.25
$ws.Range("I107").Value = 671.8889
$ws.Range("J107").Value = 365.33334
$ws.Range("K107").Value = 671.8889
$ws.Range("L107").Value = 365.33334
$ws.Range("M107").Value = 1248.1111
$ws.Range("N107").Value = -4205.33334

$ws.Range("H113").Value = 1158.5
$ws.Range("I113").Value = 1179.2858
$ws.Range("K113").Value = 1179.2858
$ws.Range("M113").Value = 990.7141999999999

$ws.Range("H136").Value = 12142
$ws.Range("I136").Value = 3999
$ws.Range("K136").Value = 11997
$ws.Range("M136").Value = -9447

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 25.666666
$ws.Range("I7").Value = 28
$ws.Range("K7").Value = 84
$ws.Range("M7").Value = 28

$ws.Range("H50").Value = 12.5
$ws.Range("I50").Value = 5
$ws.Range("J50").Value = 35
$ws.Range("K50").Value = 15
$ws.Range("L50").Value = 105
$ws.Range("M50").Value = 466
$ws.Range("N50").Value = -1067

$ws.Range("H53").Value = 12.5
$ws.Range("I53").Value = 5
$ws.Range("J53").Value = 35
$ws.Range("K53").Value = 15
$ws.Range("L53").Value = 105
$ws.Range("M53").Value = 466
$ws.Range("N53").Value = -1067

$ws.Range("H92").Value = 296.25
$ws.Range("I92").Value = 296.25
$ws.Range("K92").Value = 888.75
$ws.Range("M92").Value = 359.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40518

$ws.Range("H105").Value = 20344.5
$ws.Range("J105").Value = 20344.5
$ws.Range("L105").Value = 20344.5
$ws.Range("N105").Value = -27332.5

$ws.Range("H132").Value = 131628.5
$ws.Range("I132").Value = 185094
$ws.Range("K132").Value = 555282
$ws.Range("M132").Value = -552752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H100").Value = 2760.5
$ws.Range("I100").Value = 2760.5
$ws.Range("K100").Value = 2760.5
$ws.Range("M100").Value = -2219.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H107").Value = 531.25
$ws.Range("I107").Value = 531.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1593.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 326.25
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 728.1667
$ws.Range("J113").Value = 543.6667
$ws.Range("L113").Value = 1631.0001
$ws.Range("N113").Value = -5971.0001

$ws.Range("H136").Value = 5723.6
$ws.Range("I136").Value = 2278.2
$ws.Range("K136").Value = 6834.599999999999
$ws.Range("M136").Value = -4284.599999999999
